$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to write a value as TEXT (matching the original inlineStr/text
# cell type used throughout this worksheet) rather than letting Excel
# auto-convert the numeric-looking string into a numeric cell.
function Set-TextCellValue {
    param(
        [string]$CellRef,
        [string]$TextValue
    )

    $rng = $ws.Range($CellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $TextValue
    # Reset to the default "Normal" style so we don't leave a stray
    # number-format style applied to the cell.
    $rng.Style = "Normal"
}

# Row 32: Grand Est / Entrepreneur individuel
Set-TextCellValue "C32" "460"
Set-TextCellValue "D32" "1430342.48"

# Row 34: Grand Est / SARL
Set-TextCellValue "C34" "883"
Set-TextCellValue "D34" "6886651.94"

# Row 59: Ile-de-France / SARL
Set-TextCellValue "C59" "6820"
Set-TextCellValue "D59" "34966641.31"

# Row 82: Nouvelle-Aquitaine / SARL
Set-TextCellValue "C82" "1283"
Set-TextCellValue "D82" "10292750.17"

# Row 83: Nouvelle-Aquitaine / Societe par actions simplifiee
Set-TextCellValue "C83" "677"
Set-TextCellValue "D83" "4769396.04"
